# Apply "Added calculation of guaranteed power based on exceedance probability"

$wb = $excel.ActiveWorkbook
$wsGeneral = $wb.Worksheets.Item("General parameters")
$wsHydro   = $wb.Worksheets.Item("Hydropower plant parameters")
$wsSim     = $wb.Worksheets.Item("Simulation accuracy")

# --- Hydropower plant parameters sheet ---

# Clear the now-unused wrap style on B3:B9 (reset to Normal/default style)
$wsHydro.Range("B3:B9").Style = "Normal"

# Row 24 (f_reg): clear the C24/D24 values so the row only has A/B text
$wsHydro.Range("C24:D24").ClearContents()

# Insert a new row before row 31 for the new "p_exceedance" parameter
$wsHydro.Rows.Item(31).Insert()

$wsHydro.Cells.Item(31, 2).Value = "this percentile is used to calculate the exceedance probability of delivered power (guaranteed capacity, MW). For P90, use 90; for P95, use 95, etc."
$wsHydro.Cells.Item(31, 1).Value = "p_exceedance"
$wsHydro.Cells.Item(31, 3).Value = 90
$wsHydro.Cells.Item(31, 4).Value = 90

# Selections left by the editor
$wsGeneral.Range("C8").Select()
$wsHydro.Range("B31").Select()
$wsSim.Range("B1").Select()
